$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels (PERTANYAAN -> QUESTION, JAWABAN -> ANSWER)
$ws.Range("B1").Value = "QUESTION"
$ws.Range("C1").Value = "ANSWER"

# Move active selection to C1 (was C12)
$ws.Range("C1").Select()
